$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '2026-02-13 23:18:51'
$ws.Range('E2').NumberFormat = 'general'
$ws.Range('I2').NumberFormat = '@'
$ws.Range('I2').Value = '5.0 mm'
$ws.Range('I2').NumberFormat = 'general'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '2026-02-13 23:18:54'
$ws.Range('E3').NumberFormat = 'general'
$ws.Range('G3').NumberFormat = '@'
$ws.Range('G3').Value = '184 cm'
$ws.Range('G3').NumberFormat = 'general'
$ws.Range('H3').NumberFormat = '@'
$ws.Range('H3').Value = '89%'
$ws.Range('H3').NumberFormat = 'general'
$ws.Range('I3').NumberFormat = '@'
$ws.Range('I3').Value = '8.1 mm'
$ws.Range('I3').NumberFormat = 'general'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '2026-02-13 23:18:57'
$ws.Range('E4').NumberFormat = 'general'
$ws.Range('J4').NumberFormat = '@'
$ws.Range('J4').Value = '993.2 hPa'
$ws.Range('J4').NumberFormat = 'general'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '2026-02-13 23:18:59'
$ws.Range('E5').NumberFormat = 'general'
$ws.Range('I5').NumberFormat = '@'
$ws.Range('I5').Value = '4.5 mm'
$ws.Range('I5').NumberFormat = 'general'
$ws.Range('N5').NumberFormat = '@'
$ws.Range('N5').Value = '-4.4 °C 22:59 TU'
$ws.Range('N5').NumberFormat = 'general'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '2026-02-13 23:19:02'
$ws.Range('E6').NumberFormat = 'general'
$ws.Range('H6').NumberFormat = '@'
$ws.Range('H6').Value = '80%'
$ws.Range('H6').NumberFormat = 'general'
$ws.Range('J6').NumberFormat = '@'
$ws.Range('J6').Value = '993.2 hPa'
$ws.Range('J6').NumberFormat = 'general'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '2026-02-13 23:19:05'
$ws.Range('E7').NumberFormat = 'general'
$ws.Range('J7').NumberFormat = '@'
$ws.Range('J7').Value = '993.6 hPa'
$ws.Range('J7').NumberFormat = 'general'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '2026-02-13 23:19:08'
$ws.Range('E8').NumberFormat = 'general'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '82%'
$ws.Range('H8').NumberFormat = 'general'
$ws.Range('J8').NumberFormat = '@'
$ws.Range('J8').Value = '993.5 hPa'
$ws.Range('J8').NumberFormat = 'general'
$ws.Range('L8').NumberFormat = '@'
$ws.Range('L8').Value = '62.3 km/h - 331º 22:47 TU'
$ws.Range('L8').NumberFormat = 'general'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '2026-02-13 23:19:10'
$ws.Range('E9').NumberFormat = 'general'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '2026-02-13 23:19:13'
$ws.Range('E10').NumberFormat = 'general'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '2026-02-13 23:19:16'
$ws.Range('E11').NumberFormat = 'general'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '2026-02-13 23:19:19'
$ws.Range('E12').NumberFormat = 'general'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '2026-02-13 23:19:21'
$ws.Range('E13').NumberFormat = 'general'
$ws.Range('H13').NumberFormat = '@'
$ws.Range('H13').Value = '93%'
$ws.Range('H13').NumberFormat = 'general'
$ws.Range('J13').NumberFormat = '@'
$ws.Range('J13').Value = '996.2 hPa'
$ws.Range('J13').NumberFormat = 'general'
$ws.Range('L13').NumberFormat = '@'
$ws.Range('L13').Value = '26.6 km/h - 331º 22:20 TU'
$ws.Range('L13').NumberFormat = 'general'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '2026-02-13 23:19:24'
$ws.Range('E14').NumberFormat = 'general'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '2026-02-13 23:19:26'
$ws.Range('E15').NumberFormat = 'general'
$ws.Range('H15').NumberFormat = '@'
$ws.Range('H15').Value = '77%'
$ws.Range('H15').NumberFormat = 'general'
$ws.Range('I15').NumberFormat = '@'
$ws.Range('I15').Value = '6.4 mm'
$ws.Range('I15').NumberFormat = 'general'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '2026-02-13 23:19:29'
$ws.Range('E16').NumberFormat = 'general'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '2026-02-13 23:19:32'
$ws.Range('E17').NumberFormat = 'general'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '2026-02-13 23:19:34'
$ws.Range('E18').NumberFormat = 'general'
$ws.Range('J18').NumberFormat = '@'
$ws.Range('J18').Value = '993.4 hPa'
$ws.Range('J18').NumberFormat = 'general'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '2026-02-13 23:19:37'
$ws.Range('E19').NumberFormat = 'general'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '2026-02-13 23:19:38'
$ws.Range('E20').NumberFormat = 'general'
$ws.Range('I20').NumberFormat = '@'
$ws.Range('I20').Value = '24.6 mm'
$ws.Range('I20').NumberFormat = 'general'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '2026-02-13 23:19:40'
$ws.Range('E21').NumberFormat = 'general'
$ws.Range('H21').NumberFormat = '@'
$ws.Range('H21').Value = '93%'
$ws.Range('H21').NumberFormat = 'general'
$ws.Range('J21').NumberFormat = '@'
$ws.Range('J21').Value = '996.4 hPa'
$ws.Range('J21').NumberFormat = 'general'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '2026-02-13 23:19:41'
$ws.Range('E22').NumberFormat = 'general'
$ws.Range('L22').NumberFormat = '@'
$ws.Range('L22').Value = '76.3 km/h - 323º 22:47 TU'
$ws.Range('L22').NumberFormat = 'general'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '2026-02-13 23:19:42'
$ws.Range('E23').NumberFormat = 'general'
$ws.Range('I23').NumberFormat = '@'
$ws.Range('I23').Value = '14.2 mm'
$ws.Range('I23').NumberFormat = 'general'
$ws.Range('O23').NumberFormat = '@'
$ws.Range('O23').Value = '-4.1 °C'
$ws.Range('O23').NumberFormat = 'general'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '2026-02-13 23:19:45'
$ws.Range('E24').NumberFormat = 'general'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '2026-02-13 23:19:47'
$ws.Range('E25').NumberFormat = 'general'
$ws.Range('I25').NumberFormat = '@'
$ws.Range('I25').Value = '10.7 mm'
$ws.Range('I25').NumberFormat = 'general'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '2026-02-13 23:19:49'
$ws.Range('E26').NumberFormat = 'general'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '2026-02-13 23:19:51'
$ws.Range('E27').NumberFormat = 'general'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '2026-02-13 23:19:54'
$ws.Range('E28').NumberFormat = 'general'
$ws.Range('J28').NumberFormat = '@'
$ws.Range('J28').Value = '993.7 hPa'
$ws.Range('J28').NumberFormat = 'general'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '2026-02-13 23:19:56'
$ws.Range('E29').NumberFormat = 'general'
$ws.Range('O29').NumberFormat = '@'
$ws.Range('O29').Value = '10.9 °C'
$ws.Range('O29').NumberFormat = 'general'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '2026-02-13 23:19:59'
$ws.Range('E30').NumberFormat = 'general'
$ws.Range('J30').NumberFormat = '@'
$ws.Range('J30').Value = '993.1 hPa'
$ws.Range('J30').NumberFormat = 'general'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '2026-02-13 23:20:02'
$ws.Range('E31').NumberFormat = 'general'
$ws.Range('H31').NumberFormat = '@'
$ws.Range('H31').Value = '76%'
$ws.Range('H31').NumberFormat = 'general'
$ws.Range('I31').NumberFormat = '@'
$ws.Range('I31').Value = '5.5 mm'
$ws.Range('I31').NumberFormat = 'general'
$ws.Range('J31').NumberFormat = '@'
$ws.Range('J31').Value = '992.1 hPa'
$ws.Range('J31').NumberFormat = 'general'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '2026-02-13 23:20:05'
$ws.Range('E32').NumberFormat = 'general'
$ws.Range('O32').NumberFormat = '@'
$ws.Range('O32').Value = '4.9 °C'
$ws.Range('O32').NumberFormat = 'general'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '2026-02-13 23:20:08'
$ws.Range('E33').NumberFormat = 'general'
$ws.Range('J33').NumberFormat = '@'
$ws.Range('J33').Value = '995.2 hPa'
$ws.Range('J33').NumberFormat = 'general'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '2026-02-13 23:20:11'
$ws.Range('E34').NumberFormat = 'general'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '2026-02-13 23:20:13'
$ws.Range('E35').NumberFormat = 'general'
$ws.Range('L35').NumberFormat = '@'
$ws.Range('L35').Value = '86.4 km/h - 248º 22:35 TU'
$ws.Range('L35').NumberFormat = 'general'
$ws.Range('O35').NumberFormat = '@'
$ws.Range('O35').Value = '5.7 °C'
$ws.Range('O35').NumberFormat = 'general'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '2026-02-13 23:20:16'
$ws.Range('E36').NumberFormat = 'general'
$ws.Range('H36').NumberFormat = '@'
$ws.Range('H36').Value = '77%'
$ws.Range('H36').NumberFormat = 'general'
$ws.Range('J36').NumberFormat = '@'
$ws.Range('J36').Value = '993.2 hPa'
$ws.Range('J36').NumberFormat = 'general'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '2026-02-13 23:20:18'
$ws.Range('E37').NumberFormat = 'general'
$ws.Range('J37').NumberFormat = '@'
$ws.Range('J37').Value = '995.1 hPa'
$ws.Range('J37').NumberFormat = 'general'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '2026-02-13 23:20:21'
$ws.Range('E38').NumberFormat = 'general'
$ws.Range('H38').NumberFormat = '@'
$ws.Range('H38').Value = '80%'
$ws.Range('H38').NumberFormat = 'general'
$ws.Range('N38').NumberFormat = '@'
$ws.Range('N38').Value = '7.7 °C 22:46 TU'
$ws.Range('N38').NumberFormat = 'general'
$ws.Range('O38').NumberFormat = '@'
$ws.Range('O38').Value = '9.5 °C'
$ws.Range('O38').NumberFormat = 'general'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '2026-02-13 23:20:23'
$ws.Range('E39').NumberFormat = 'general'
$ws.Range('I39').NumberFormat = '@'
$ws.Range('I39').Value = '20.1 mm'
$ws.Range('I39').NumberFormat = 'general'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '2026-02-13 23:20:26'
$ws.Range('E40').NumberFormat = 'general'
$ws.Range('I40').NumberFormat = '@'
$ws.Range('I40').Value = '17.3 mm'
$ws.Range('I40').NumberFormat = 'general'
$ws.Range('J40').NumberFormat = '@'
$ws.Range('J40').Value = '996.9 hPa'
$ws.Range('J40').NumberFormat = 'general'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '2026-02-13 23:20:29'
$ws.Range('E41').NumberFormat = 'general'
$ws.Range('J41').NumberFormat = '@'
$ws.Range('J41').Value = '994.0 hPa'
$ws.Range('J41').NumberFormat = 'general'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '2026-02-13 23:20:32'
$ws.Range('E42').NumberFormat = 'general'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '2026-02-13 23:20:35'
$ws.Range('E43').NumberFormat = 'general'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '2026-02-13 23:20:37'
$ws.Range('E44').NumberFormat = 'general'
$ws.Range('H44').NumberFormat = '@'
$ws.Range('H44').Value = '91%'
$ws.Range('H44').NumberFormat = 'general'
$ws.Range('I44').NumberFormat = '@'
$ws.Range('I44').Value = '11.4 mm'
$ws.Range('I44').NumberFormat = 'general'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '2026-02-13 23:20:40'
$ws.Range('E45').NumberFormat = 'general'
$ws.Range('I45').NumberFormat = '@'
$ws.Range('I45').Value = '3.0 mm'
$ws.Range('I45').NumberFormat = 'general'
$ws.Range('J45').NumberFormat = '@'
$ws.Range('J45').Value = '992.9 hPa'
$ws.Range('J45').NumberFormat = 'general'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '2026-02-13 23:20:43'
$ws.Range('E46').NumberFormat = 'general'
$ws.Range('H46').NumberFormat = '@'
$ws.Range('H46').Value = '86%'
$ws.Range('H46').NumberFormat = 'general'
